# Applies yellow highlighting to specific key terms throughout the
# requirements document, splitting runs as needed so only the target
# word(s) carry the <w:highlight w:val="yellow"/> run property.

$d = $word.ActiveDocument

function Highlight-Word {
    param(
        [string]$Phrase,   # a phrase that is unique in the document
        [string]$Word      # the substring inside $Phrase to highlight (unique within $Phrase)
    )

    $outer = $d.Content
    $outer.Find.Execute($Phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

    $inner = $d.Range($outer.Start, $outer.End)
    $inner.Find.Execute($Word, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

    $inner.Font.HighlightColorIndex = 7
}

# 1. "...tres tipos de usuarios: Alumnos que"
Highlight-Word "La base de datos será utilizada por tres tipos de usuarios: Alumnos que" "usuarios"

# 2. "De los alumnos se desea almacenar: "
Highlight-Word "De los alumnos se desea almacenar: " "alumnos"

# 3. "De las carreras se "
Highlight-Word "De las carreras se " "carreras"

# 4. "Para los docentes es necesario almacenar: "
Highlight-Word "Para los docentes es necesario almacenar: " "docentes"

# 5. "...distintas ediciones del concurso"
Highlight-Word "Se requiere almacenar información de las distintas ediciones del concurso" "ediciones"

# 6. "De los equipos participantes en cada una de las ediciones del concurso se"
Highlight-Word "De los equipos participantes en cada una de las ediciones del concurso se" "equipos"

# 7. "De las categorías de los equipos se desea almacenar: un identificador único,"
Highlight-Word "De las categorías de los equipos se desea almacenar: un identificador único," "categorías"

# 8. Existing standalone runs "banco" / " " / "de" / " " / "problemas" all get highlighted
#    (they already form separate runs, so this just adds the property).
Highlight-Word "banco de problemas para las distintas ediciones del concurso" "banco de problemas"

# 9. Existing standalone run "edición" (inside "...para cada edición del concurso se...")
Highlight-Word "para cada edición del concurso se seleccionen" "edición"

# 10. "También se necesita almacenar los problemas resueltos por cada equipo"
Highlight-Word "También se necesita almacenar los problemas resueltos por cada equipo" "problemas resueltos"

Write-Output "Done"
